$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text/"@" number format on every cell we are about to write so that
# values such as "1.00", "7.71", "0.0830" are preserved exactly as text and
# are not re-interpreted/rounded by Excel as floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.315.76'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.232.64'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.91'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +6.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.52'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.582'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +7.75%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.562'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +7.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.18'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +9.75%  '
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0830'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.26%  '
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.71'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +11.49%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.42%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.865'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.568.93'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.41'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +9.28%  '
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.227.97'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.184.18'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.33%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.32'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +12.75%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0970'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.97%  '
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.60'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.71%  '
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'PancakeSwap'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.21'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +6.32%  '
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.25'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.93%  '
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.49'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.76%  '
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.17'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +9.90%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.59%  '
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.06'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.40%  '
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.94%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +8.90%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.43'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.26%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.64'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +16.76%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.28'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.85%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0869'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.44%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '156.84'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.64%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.68'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.26'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +23.45%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.121'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.14%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.86'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +12.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.42'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +11.58%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.106'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +24.31%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.70'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +13.14%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0321'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.18%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.31'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +33.46%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.782.74'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +8.80%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.203'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.46%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.27'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.96'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.61%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.29'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.38%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.15'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +6.10%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.38'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.43%  '
